{"js": "// 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", \"Replace\");\n}\n\n// 2) Split the mailing-address paragraph (\"919 Story Road, San Jose CA 95122\")\n//    into two separate paragraphs: \"919 Story Road\" and \"San Jose, CA 95122\".\n//    Only the first (stand-alone) occurrence is affected -- the identical\n//    text that lives inside the account-info table is left untouched.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"919 Story Road, San Jose CA 95122\";\nconst candidates = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    candidates.push(paragraphs.items[i]);\n  }\n}\n// Confirm which candidate (if more than one) is NOT inside a table.\nconst candidateParents = candidates.map((cand) => cand.parentTableOrNullObject);\ncandidateParents.forEach((pt) => pt.load(\"isNullObject\"));\nawait context.sync();\nlet addressPara = null;\nfor (let i = 0; i < candidates.length; i++) {\n  if (candidateParents[i].isNullObject) {\n    addressPara = candidates[i];\n    break;\n  }\n}\n\nif (addressPara) {\n  addressPara.insertParagraph(\"San Jose, CA 95122\", \"After\");\n  addressPara.getRange().insertText(\"919 Story Road\", \"Replace\");\n}\nawait context.sync();\n\n// 3) Remove the now-superfluous blank \"No Spacing\" paragraph that used to sit\n//    right after \"...Board of Directors\".\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items/text,items/style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (\n    paragraphs2.items[i].text === \"\" &&\n    paragraphs2.items[i].style === \"No Spacing\" &&\n    i > 0 &&\n    paragraphs2.items[i - 1].text.indexOf(\"Board of Directors\") !== -1\n  ) {\n    paragraphs2.items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"September 19, 2025\", $false, $false, $false, $false, $false, $true, 1, $false, \"September 21, 2025\", 2) | Out-Null\n\n# 2) Split the stand-alone mailing-address paragraph\n#    (\"919 Story Road, San Jose CA 95122\") into two paragraphs:\n#    \"919 Story Road\" and \"San Jose, CA 95122\". The identical text that\n#    lives inside the account-info table must be left untouched, so we\n#    explicitly skip any paragraph that is inside a table.\n$oldAddress = \"919 Story Road, San Jose CA 95122\"\n$addrPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq $oldAddress -and -not $p.Range.Information(12)) {\n        $addrPara = $p\n        break\n    }\n}\n\nif ($addrPara -ne $null) {\n    $r = $addrPara.Range\n    $bodyRange = $d.Range($r.Start, $r.End - 1)\n    $bodyRange.Text = \"919 Story Road\" + [char]13 + \"San Jose, CA 95122\"\n}\n\n# Re-fetch the freshly created second paragraph by its text (robust against\n# any index shifting) and make sure its run formatting (Arial, 11pt, incl.\n# complex-script variants) matches the rest of the letter.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"San Jose, CA 95122\") {\n        $r3 = $p.Range\n        $r3.Font.Name = \"Arial\"\n        $r3.Font.NameBi = \"Arial\"\n        $r3.Font.Size = 11\n        $r3.Font.SizeBi = 11\n        break\n    }\n}\n\n# 3) Remove the now-superfluous blank \"No Spacing\" paragraph that used to\n#    sit right after \"...Board of Directors\".\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"\" -and $p.Style.NameLocal -eq \"No Spacing\" -and $i -gt 1) {\n        $prev = $d.Paragraphs.Item($i - 1)\n        $prevTxt = $prev.Range.Text.TrimEnd([char]13, [char]7)\n        if ($prevTxt -like \"*Board of Directors*\") {\n            $p.Range.Delete()\n            break\n        }\n    }\n}\n"}
